# Update column G ("K") values on the data sheet.
# These new values come from a regenerated save_data pipeline (per commit
# message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals"). Row 34 is unchanged (already 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 5
    6  = 2
    7  = 3
    8  = 2
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 2
    15 = 1
    16 = 1
    17 = 2
    18 = 3
    19 = 2
    20 = 0
    21 = 1
    22 = 0
    23 = 2
    24 = 0
    25 = 0
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    30 = 1
    31 = 0
    32 = 2
    33 = 2
    35 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
